# Benchmark update: new Chroma numbers from "anton input".
#
# Column D ("Qdrant" in the original sheet) shifts to column E, and a brand
# new "Chroma" column of data is inserted in its place at column D. The new
# D1 header ("Chroma") is emphasised with bold text, a thin box border and
# center/top alignment (matching the look already used for the other header
# cells); E1 keeps the normal header look.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ------------------------------------------------
# D1 becomes "Chroma" (previously the old E1 label), E1 becomes "Qdrant"
# (previously the old D1 label) - i.e. the two header labels swap places.
$ws.Range("D1").Value = "Chroma"
$ws.Range("E1").Value = "Qdrant"

# Emphasize the new D1 header cell: bold font, thin border all around,
# centered horizontally and aligned to the top vertically.
$d1 = $ws.Range("D1")
$d1.Font.Bold = $true
$d1.Borders.LineStyle = 1
$d1.Borders.Weight = 2
$d1.HorizontalAlignment = -4108
$d1.VerticalAlignment = -4160

# --- Data rows (rows 2-5) ----------------------------------------------
# The values that used to live in column D (old "Qdrant" numbers) move over
# to column E (now labeled "Qdrant"). Column D receives brand new "Chroma"
# benchmark numbers.
$ws.Range("E2").Value = 0.04868075083009897747
$ws.Range("E3").Value = 0.04485013209050522226
$ws.Range("E4").Value = 0.37011745250085370529
$ws.Range("E5").Value = 0.07904000750044361256

$ws.Range("D2").Value = 0.01615460458211601114
$ws.Range("D3").Value = 0.01155560291837900988
$ws.Range("D4").Value = 0.14928332374896849388
$ws.Range("D5").Value = 0.00958641249919310254

# --- Selection -----------------------------------------------------------
[void]$ws.Range("L9").Select()
